$wb = $excel.ActiveWorkbook

# --- Sheet "missing_values" (sheet1) ---
$ws1 = $wb.Worksheets.Item("missing_values")

# Rows with mirrored D:I columns (B/C values duplicated into D/E, F/G, H/I pairs)
$ws1.Range("B4").Value = 176
$ws1.Range("C4").Value = 1.3573962671602653
$ws1.Range("D4").Value = 176
$ws1.Range("E4").Value = 1.3573962671602653
$ws1.Range("F4").Value = 176
$ws1.Range("G4").Value = 1.3573962671602653
$ws1.Range("H4").Value = 176
$ws1.Range("I4").Value = 1.3573962671602653
$ws1.Range("B5").Value = 651
$ws1.Range("C5").Value = 5.0208236927348455
$ws1.Range("D5").Value = 651
$ws1.Range("E5").Value = 5.0208236927348455
$ws1.Range("F5").Value = 651
$ws1.Range("G5").Value = 5.0208236927348455
$ws1.Range("H5").Value = 651
$ws1.Range("I5").Value = 5.0208236927348455
$ws1.Range("B6").Value = 1687
$ws1.Range("C6").Value = 13.010951719882771
$ws1.Range("D6").Value = 1687
$ws1.Range("E6").Value = 13.010951719882771
$ws1.Range("F6").Value = 1687
$ws1.Range("G6").Value = 13.010951719882771
$ws1.Range("H6").Value = 1687
$ws1.Range("I6").Value = 13.010951719882771
$ws1.Range("B7").Value = 41
$ws1.Range("C7").Value = 0.31621163041801637
$ws1.Range("D7").Value = 41
$ws1.Range("E7").Value = 0.31621163041801637
$ws1.Range("F7").Value = 41
$ws1.Range("G7").Value = 0.31621163041801637
$ws1.Range("H7").Value = 41
$ws1.Range("I7").Value = 0.31621163041801637
$ws1.Range("B8").Value = 2379
$ws1.Range("C8").Value = 18.347987043035634
$ws1.Range("D8").Value = 2379
$ws1.Range("E8").Value = 18.347987043035634
$ws1.Range("F8").Value = 2379
$ws1.Range("G8").Value = 18.347987043035634
$ws1.Range("H8").Value = 2379
$ws1.Range("I8").Value = 18.347987043035634
$ws1.Range("B9").Value = 10377
$ws1.Range("C9").Value = 80.032392410920878
$ws1.Range("D9").Value = 10377
$ws1.Range("E9").Value = 80.032392410920878
$ws1.Range("F9").Value = 10377
$ws1.Range("G9").Value = 80.032392410920878
$ws1.Range("H9").Value = 10377
$ws1.Range("I9").Value = 80.032392410920878
$ws1.Range("B10").Value = 12966
$ws1.Range("C10").Value = 100
$ws1.Range("D10").Value = 12966
$ws1.Range("E10").Value = 100
$ws1.Range("F10").Value = 12966
$ws1.Range("G10").Value = 100
$ws1.Range("H10").Value = 12966
$ws1.Range("I10").Value = 100
$ws1.Range("B15").Value = 23447
$ws1.Range("C15").Value = 374.6724193032918
$ws1.Range("D15").Value = 23447
$ws1.Range("E15").Value = 374.6724193032918
$ws1.Range("F15").Value = 23447
$ws1.Range("G15").Value = 374.6724193032918
$ws1.Range("H15").Value = 23447
$ws1.Range("I15").Value = 374.6724193032918
$ws1.Range("B16").Value = 147
$ws1.Range("C16").Value = 2.348993288590604
$ws1.Range("D16").Value = 147
$ws1.Range("E16").Value = 2.348993288590604
$ws1.Range("F16").Value = 147
$ws1.Range("G16").Value = 2.348993288590604
$ws1.Range("H16").Value = 147
$ws1.Range("I16").Value = 2.348993288590604
$ws1.Range("B17").Value = 18
$ws1.Range("C17").Value = 0.28763183125599234
$ws1.Range("D17").Value = 18
$ws1.Range("E17").Value = 0.28763183125599234
$ws1.Range("F17").Value = 18
$ws1.Range("G17").Value = 0.28763183125599234
$ws1.Range("H17").Value = 18
$ws1.Range("I17").Value = 0.28763183125599234
$ws1.Range("B18").Value = 6
$ws1.Range("C18").Value = 0.095877277085330767
$ws1.Range("D18").Value = 6
$ws1.Range("E18").Value = 0.095877277085330767
$ws1.Range("F18").Value = 6
$ws1.Range("G18").Value = 0.095877277085330767
$ws1.Range("H18").Value = 6
$ws1.Range("I18").Value = 0.095877277085330767
$ws1.Range("B19").Value = 171
$ws1.Range("C19").Value = 2.7325023969319271
$ws1.Range("D19").Value = 171
$ws1.Range("E19").Value = 2.7325023969319271
$ws1.Range("F19").Value = 171
$ws1.Range("G19").Value = 2.7325023969319271
$ws1.Range("H19").Value = 171
$ws1.Range("I19").Value = 2.7325023969319271
$ws1.Range("B20").Value = 6024
$ws1.Range("C20").Value = 96.260786193672104
$ws1.Range("D20").Value = 6024
$ws1.Range("E20").Value = 96.260786193672104
$ws1.Range("F20").Value = 6024
$ws1.Range("G20").Value = 96.260786193672104
$ws1.Range("H20").Value = 6024
$ws1.Range("I20").Value = 96.260786193672104
$ws1.Range("B21").Value = 6258
$ws1.Range("C21").Value = 100
$ws1.Range("D21").Value = 6258
$ws1.Range("E21").Value = 100
$ws1.Range("F21").Value = 6258
$ws1.Range("G21").Value = 100
$ws1.Range("H21").Value = 6258
$ws1.Range("I21").Value = 100
$ws1.Range("B26").Value = 15
$ws1.Range("C26").Value = 1.6816143497757847
$ws1.Range("D26").Value = 15
$ws1.Range("E26").Value = 1.6816143497757847
$ws1.Range("F26").Value = 15
$ws1.Range("G26").Value = 1.6816143497757847
$ws1.Range("H26").Value = 15
$ws1.Range("I26").Value = 1.6816143497757847
$ws1.Range("B27").Value = 59
$ws1.Range("C27").Value = 6.6143497757847527
$ws1.Range("D27").Value = 59
$ws1.Range("E27").Value = 6.6143497757847527
$ws1.Range("F27").Value = 59
$ws1.Range("G27").Value = 6.6143497757847527
$ws1.Range("H27").Value = 59
$ws1.Range("I27").Value = 6.6143497757847527
$ws1.Range("B28").Value = 2
$ws1.Range("C28").Value = 0.22421524663677131
$ws1.Range("D28").Value = 2
$ws1.Range("E28").Value = 0.22421524663677131
$ws1.Range("F28").Value = 2
$ws1.Range("G28").Value = 0.22421524663677131
$ws1.Range("H28").Value = 2
$ws1.Range("I28").Value = 0.22421524663677131
$ws1.Range("B29").Value = 61
$ws1.Range("C29").Value = 6.8385650224215251
$ws1.Range("D29").Value = 61
$ws1.Range("E29").Value = 6.8385650224215251
$ws1.Range("F29").Value = 61
$ws1.Range("G29").Value = 6.8385650224215251
$ws1.Range("H29").Value = 61
$ws1.Range("I29").Value = 6.8385650224215251
$ws1.Range("B30").Value = 816
$ws1.Range("C30").Value = 91.479820627802695
$ws1.Range("D30").Value = 816
$ws1.Range("E30").Value = 91.479820627802695
$ws1.Range("F30").Value = 816
$ws1.Range("G30").Value = 91.479820627802695
$ws1.Range("H30").Value = 816
$ws1.Range("I30").Value = 91.479820627802695
$ws1.Range("B31").Value = 892
$ws1.Range("C31").Value = 100
$ws1.Range("D31").Value = 892
$ws1.Range("E31").Value = 100
$ws1.Range("F31").Value = 892
$ws1.Range("G31").Value = 100
$ws1.Range("H31").Value = 892
$ws1.Range("I31").Value = 100

# Rows 37-40: only B/C updated (no D:I mirrors)
$ws1.Range("B37").Value = 176
$ws1.Range("C37").Value = 2.1999999999999997
$ws1.Range("B38").Value = 64
$ws1.Range("C38").Value = 0.80000000000000004
$ws1.Range("B39").Value = 176
$ws1.Range("C39").Value = 2.1999999999999997
$ws1.Range("B40").Value = 7753
$ws1.Range("C40").Value = 96.912499999999994

# --- Sheet "profile_missing_values" (sheet2) ---
$ws2 = $wb.Worksheets.Item("profile_missing_values")

# Clear header cells B2/C2 (previously held "Labor income%"/"Pensions%" labels)
$ws2.Range("B2").Value = ""
$ws2.Range("C2").Value = ""

$ws2.Range("B3").Value = 33.039092055485497
$ws2.Range("C3").Value = 54.970760233918128
$ws2.Range("B4").Value = 66.960907944514503
$ws2.Range("C4").Value = 45.029239766081872
$ws2.Range("B6").Value = 0.33627574611181166
$ws2.Range("C6").Value = 0.58479532163742687
$ws2.Range("B7").Value = 12.35813366960908
$ws2.Range("C7").Value = 0
$ws2.Range("B8").Value = 22.320302648171499
$ws2.Range("C8").Value = 0
$ws2.Range("B9").Value = 23.245060949978981
$ws2.Range("C9").Value = 1.1695906432748537
$ws2.Range("B10").Value = 20.638923917612441
$ws2.Range("C10").Value = 5.2631578947368416
$ws2.Range("B11").Value = 14.543926019335856
$ws2.Range("C11").Value = 39.1812865497076
$ws2.Range("B12").Value = 6.557377049180328
$ws2.Range("C12").Value = 53.801169590643269
$ws2.Range("B14").Value = 24.506094997898277
$ws2.Range("C14").Value = 38.011695906432749
$ws2.Range("B15").Value = 32.492643968053805
$ws2.Range("C15").Value = 15.789473684210526
$ws2.Range("B16").Value = 29.087852038671713
$ws2.Range("C16").Value = 9.9415204678362574
$ws2.Range("B17").Value = 10.928961748633879
$ws2.Range("C17").Value = 12.865497076023392
$ws2.Range("B18").Value = 2.8583438419503993
$ws2.Range("C18").Value = 22.807017543859647
$ws2.Range("B19").Value = 0.12610340479192939
$ws2.Range("C19").Value = 0.58479532163742687
$ws2.Range("B21").Value = 0.29424127784783521
$ws2.Range("C21").Value = 0
$ws2.Range("B22").Value = 0.12610340479192939
$ws2.Range("C22").Value = 0.58479532163742687
$ws2.Range("B23").Value = 26.691887347625055
$ws2.Range("C23").Value = 49.122807017543856
$ws2.Range("B24").Value = 46.784363177805801
$ws2.Range("C24").Value = 23.391812865497073
$ws2.Range("B25").Value = 5.9688944934846573
$ws2.Range("C25").Value = 8.1871345029239766
$ws2.Range("B26").Value = 14.754098360655737
$ws2.Range("C26").Value = 8.7719298245614024
$ws2.Range("B27").Value = 0.67255149222362332
$ws2.Range("C27").Value = 0
$ws2.Range("B28").Value = 4.7078604455653634
$ws2.Range("C28").Value = 9.9415204678362574
$ws2.Range("B30").Value = 9.3316519546027745
$ws2.Range("C30").Value = 8.7719298245614024
$ws2.Range("B31").Value = 6.9356872635561162
$ws2.Range("C31").Value = 2.9239766081871341
$ws2.Range("B32").Value = 10.004203446826397
$ws2.Range("C32").Value = 2.9239766081871341
$ws2.Range("B33").Value = 66.876839007986547
$ws2.Range("C33").Value = 12.280701754385964
$ws2.Range("B34").Value = 0.50441361916771754
$ws2.Range("C34").Value = 0
$ws2.Range("B35").Value = 3.6990332072299283
$ws2.Range("C35").Value = 0
$ws2.Range("B36").Value = 2.3118957545187055
$ws2.Range("C36").Value = 1.7543859649122806
$ws2.Range("B37").Value = 0.33627574611181166
$ws2.Range("C37").Value = 71.345029239766077
$ws2.Range("B39").Value = 18.032786885245901
$ws2.Range("C39").Value = 5.8479532163742682
$ws2.Range("B40").Value = 0.58848255569567043
$ws2.Range("C40").Value = 1.1695906432748537
$ws2.Range("B41").Value = 1.807482135350988
$ws2.Range("C41").Value = 0.58479532163742687
$ws2.Range("B42").Value = 1.3871374527112232
$ws2.Range("C42").Value = 0.58479532163742687
$ws2.Range("B43").Value = 5.6326187473728462
$ws2.Range("C43").Value = 2.3391812865497075
$ws2.Range("B44").Value = 20.176544766708702
$ws2.Range("C44").Value = 1.7543859649122806
$ws2.Range("B45").Value = 11.601513240857503
$ws2.Range("C45").Value = 2.3391812865497075
$ws2.Range("B46").Value = 3.278688524590164
$ws2.Range("C46").Value = 0.58479532163742687
$ws2.Range("B47").Value = 8.9533417402269855
$ws2.Range("C47").Value = 7.0175438596491224
$ws2.Range("B48").Value = 28.205128205128204
$ws2.Range("C48").Value = 6.4327485380116958
$ws2.Range("B49").Value = 0.33627574611181166
$ws2.Range("C49").Value = 71.345029239766077
$ws2.Range("B51").Value = 19.630096679277006
$ws2.Range("C51").Value = 10.526315789473683
$ws2.Range("B52").Value = 80.369903320722997
$ws2.Range("C52").Value = 89.473684210526315

Write-Host "Edit applied successfully"
